$p = $ppt.ActivePresentation

# --- 1) Table on slide 5: switch the table's style to the new built-in style id ---
$s5 = $p.Slides.Item(5)
for ($i = 1; $i -le $s5.Shapes.Count; $i++) {
    $shp = $s5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{8CFE391E-50FE-40CD-BD48-370C020846E4}")
    }
}

# --- 2) Re-point the deck's live theme (master/theme2.xml) back to the classic
#        "Office" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink), undoing the
#        "Integral" / "Red Violet" colors it currently carries. ---
$master = $p.Designs.Item(1).SlideMaster
$scheme = $master.ColorScheme

$officeColors = @{
    1  = 0          # dk1      000000
    2  = 16777215   # lt1      FFFFFF
    3  = 6968388     # dk2      44546A
    4  = 15132391    # lt2      E7E6E6
    5  = 13998939    # accent1  5B9BD5
    6  = 3243501     # accent2  ED7D31
    7  = 10855845    # accent3  A5A5A5
    8  = 49407       # accent4  FFC000
    9  = 12874308    # accent5  4472C4
    10 = 4697456     # accent6  70AD47
    11 = 12673797    # hlink    0563C1
    12 = 7491477     # folHlink 954F72
}

foreach ($idx in 1..12) {
    $scheme.Colors($idx).RGB = $officeColors[$idx]
}
